$wb = $excel.ActiveWorkbook

# --- Sheet "year": append a new scaling rule row ---
$wsYear = $wb.Worksheets.Item("year")
$wsYear.Cells.Item(4, 1).Value = "lux"
$wsYear.Cells.Item(4, 2).Value = "Road"
$wsYear.Cells.Item(4, 3).Value = "NA"
$wsYear.Cells.Item(4, 4).Value = "NA"
$wsYear.Cells.Item(4, 5).Value = "NA"
$wsYear.Cells.Item(4, 6).Value = 1990
$wsYear.Cells.Item(4, 7).Value = 2020
$wsYear.Cells.Item(4, 8).Value = "NA"
$wsYear.Cells.Item(4, 9).Value = "Avoid imlied Nox EF dip 1986-1989"

# --- Sheet "map": row 42 ---
# Clear the scaling_sector value in B42 (was "Other") and add a comment in E42.
$wsMap = $wb.Worksheets.Item("map")
$wsMap.Cells.Item(42, 2).ClearContents()
$wsMap.Cells.Item(42, 5).Value = "Problematic to lump this with stationary sectors, and probably something other than rail here, so try leaving out"

# --- Update selections / active views to match where the user ended up ---
[void]$wsYear.Range("A5").Select()
[void]$wsMap.Range("E43").Select()
[void]$wsMap.Activate()
